# Applies:
#  1. Column C ("Förändrad") date serial bumped from 45184 to 45186 for every
#     data row (rows 2..257).
#  2. The HYPERLINK() formulas in columns S, T, V, W, X, Y for rows 2 and 3
#     gain a second (friendly-name) argument equal to the row's "Beteckning"
#     (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 257
$oldDate = 45184
$newDate = 45186

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    if ($cCell.Value2 -eq $oldDate) {
        $cCell.Value = $newDate
    }
}

$hyperlinkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y

foreach ($r in @(2, 3)) {
    $label = $ws.Cells.Item($r, 1).Value2   # column A, e.g. "A 12589-2021"
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ($f -and $f.StartsWith("=HYPERLINK(")) {
            $inner = $f.Substring(11, $f.Length - 12)   # strip "=HYPERLINK(" and trailing ")"
            $newFormula = "=HYPERLINK(" + $inner + ", " + '"' + $label + '"' + ")"
            $cell.Formula = $newFormula
        }
    }
}
